$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 84.333336
$ws.Range("I29").Value = 84.333336
$ws.Range("K29").Value = 253.000008
$ws.Range("M29").Value = 27.99999199999999

$ws.Range("H33").Value = 626.4375
$ws.Range("I33").Value = 323.0909
$ws.Range("K33").Value = 323.0909
$ws.Range("M33").Value = -94.09089999999998

$ws.Range("H40").Value = 3999.95

$ws.Range("H43").Value = 2605.9443
$ws.Range("I43").Value = 2400.75
$ws.Range("K43").Value = 2400.75
$ws.Range("M43").Value = -2331.75

$ws.Range("H70").Value = 83334664
$ws.Range("I70").Value = 999
$ws.Range("K70").Value = 2997
$ws.Range("M70").Value = -2727

$ws.Range("H73").Value = 83334664
$ws.Range("I73").Value = 999
$ws.Range("K73").Value = 2997
$ws.Range("M73").Value = -2061

$ws.Range("H82").Value = 220
$ws.Range("I82").Value = 220
$ws.Range("K82").Value = 660
$ws.Range("M82").Value = -254

$ws.Range("H85").Value = 220
$ws.Range("I85").Value = 220
$ws.Range("K85").Value = 660
$ws.Range("M85").Value = 744

$ws.Range("H115").Value = 521.8
$ws.Range("I115").Value = 521.8
$ws.Range("K115").Value = 1565.4
$ws.Range("M115").Value = 1.600000000000136

$ws.Range("H132").Value = 2437.7932
$ws.Range("I132").Value = 2268.5
$ws.Range("J132").Value = 2714.818
$ws.Range("K132").Value = 6805.5
$ws.Range("L132").Value = 8144.454000000001
$ws.Range("M132").Value = -4275.5
$ws.Range("N132").Value = -13204.454

$ws.Range("H135").Value = 1077.7222
$ws.Range("J135").Value = 3134.1428
$ws.Range("L135").Value = 28207.2852
$ws.Range("N135").Value = -33277.2852

$ws.Range("H137").Value = 32941.668
$ws.Range("I137").Value = 38105.914
$ws.Range("K137").Value = 114317.742
$ws.Range("M137").Value = -111767.742

$ws.Range("H138").Value = 2525.81
$ws.Range("I138").Value = 1320.2727
$ws.Range("J138").Value = 2674.809
$ws.Range("K138").Value = 3960.8181
$ws.Range("L138").Value = 8024.427000000001
$ws.Range("M138").Value = 1179.1819
$ws.Range("N138").Value = -18304.427

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 8091.4287
$ws.Range("I122").Value = 5878.1816
$ws.Range("K122").Value = 17634.5448
$ws.Range("M122").Value = -15184.5448

$ws.Range("H132").Value = 3482.524
$ws.Range("I132").Value = 3008.125
$ws.Range("J132").Value = 5000.6
$ws.Range("K132").Value = 9024.375
$ws.Range("L132").Value = 15001.8
$ws.Range("M132").Value = -6494.375
$ws.Range("N132").Value = -20061.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 6198.9
$ws.Range("I128").Value = 6198.9
$ws.Range("K128").Value = 18596.7
$ws.Range("M128").Value = -16106.7

$ws.Range("H134").Value = 1589352.8
$ws.Range("I134").Value = 1756116.1
$ws.Range("K134").Value = 5268348.300000001
$ws.Range("M134").Value = -5265813.300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 6000
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 10000
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 10000
$ws.Range("M55").Value = -1685
$ws.Range("N55").Value = -10630

$ws.Range("H58").Value = 3105.566
$ws.Range("I58").Value = 2930.8809
$ws.Range("J58").Value = 3772.5454
$ws.Range("K58").Value = 2930.8809
$ws.Range("L58").Value = 3772.5454
$ws.Range("M58").Value = -2727.8809
$ws.Range("N58").Value = -4178.5454

$ws.Range("H60").Value = 9093
$ws.Range("I60").Value = 9093
$ws.Range("K60").Value = 9093
$ws.Range("M60").Value = -8582

$ws.Range("H86").Value = 46728.684
$ws.Range("J86").Value = 44545.16
$ws.Range("L86").Value = 44545.16
$ws.Range("N86").Value = -46791.16

$ws.Range("H89").Value = 46728.684
$ws.Range("J89").Value = 44545.16
$ws.Range("L89").Value = 222725.8
$ws.Range("N89").Value = -233957.8

$ws.Range("H122").Value = 4711.0527
$ws.Range("I122").Value = 3518.1538
$ws.Range("K122").Value = 10554.4614
$ws.Range("M122").Value = -8104.4614

$ws.Range("H125").Value = 70689.39999999999
$ws.Range("J125").Value = 70689.39999999999
$ws.Range("L125").Value = 70689.39999999999
$ws.Range("N125").Value = -75609.39999999999

$ws.Range("H134").Value = 2523.4443
$ws.Range("I134").Value = 2266.3333
$ws.Range("J134").Value = 2652
$ws.Range("K134").Value = 6798.999899999999
$ws.Range("L134").Value = 7956
$ws.Range("M134").Value = -4263.999899999999
$ws.Range("N134").Value = -13026

$ws.Range("H136").Value = 3105.566
$ws.Range("I136").Value = 2930.8809
$ws.Range("J136").Value = 3772.5454
$ws.Range("K136").Value = 8792.6427
$ws.Range("L136").Value = 11317.6362
$ws.Range("M136").Value = -6242.6427
$ws.Range("N136").Value = -16417.6362

$ws.Range("H141").Value = 816326.5600000001
$ws.Range("J141").Value = 816326.5600000001
$ws.Range("L141").Value = 816326.5600000001
$ws.Range("N141").Value = -826686.5600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4332.1113
$ws.Range("J113").Value = 4748.625
$ws.Range("L113").Value = 14245.875
$ws.Range("N113").Value = -18585.875

$ws.Range("H138").Value = 42601210
$ws.Range("J138").Value = 71000664
$ws.Range("L138").Value = 213001992
$ws.Range("N138").Value = -213012272

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2686.4546
$ws.Range("I126").Value = 2690.182
$ws.Range("K126").Value = 8070.545999999999
$ws.Range("M126").Value = -5600.545999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 2513500
$ws.Range("I18").Value = 27000
$ws.Range("K18").Value = 27000
$ws.Range("M18").Value = -26828

$ws.Range("H22").Value = 2588.0625
$ws.Range("I22").Value = 1835.7
$ws.Range("J22").Value = 3842
$ws.Range("K22").Value = 1835.7
$ws.Range("L22").Value = 3842
$ws.Range("M22").Value = -1540.7
$ws.Range("N22").Value = -4432

$ws.Range("H27").Value = 2588.0625
$ws.Range("I27").Value = 1835.7
$ws.Range("J27").Value = 3842
$ws.Range("K27").Value = 1835.7
$ws.Range("L27").Value = 3842
$ws.Range("M27").Value = -1728.7
$ws.Range("N27").Value = -4056

$ws.Range("H40").Value = 4000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H46").Value = 8910.409
$ws.Range("J46").Value = 9191.857
$ws.Range("L46").Value = 9191.857
$ws.Range("N46").Value = -9567.857

$ws.Range("H68").Value = 1945.2
$ws.Range("I68").Value = 1680.875
$ws.Range("J68").Value = 3002.5
$ws.Range("K68").Value = 1680.875
$ws.Range("L68").Value = 3002.5
$ws.Range("M68").Value = -931.875
$ws.Range("N68").Value = -4500.5

$ws.Range("H71").Value = 1945.2
$ws.Range("I71").Value = 1680.875
$ws.Range("J71").Value = 3002.5
$ws.Range("K71").Value = 8404.375
$ws.Range("L71").Value = 15012.5
$ws.Range("M71").Value = -4660.375
$ws.Range("N71").Value = -22500.5

$ws.Range("H131").Value = 82633
$ws.Range("J131").Value = 82633
$ws.Range("L131").Value = 82633
$ws.Range("N131").Value = -92713

$ws.Range("H136").Value = 10359
$ws.Range("I136").Value = 6698.75
$ws.Range("K136").Value = 20096.25
$ws.Range("M136").Value = -17546.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 141329.67
$ws.Range("J101").Value = 141329.67
$ws.Range("L101").Value = 141329.67
$ws.Range("N101").Value = -147819.67
